$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 15222.5
$ws.Range("I62").Value = 19300
$ws.Range("J62").Value = 2990
$ws.Range("K62").Value = 19300
$ws.Range("L62").Value = 2990
$ws.Range("M62").Value = -18676
$ws.Range("N62").Value = -4238
$ws.Range("H65").Value = 15222.5
$ws.Range("I65").Value = 19300
$ws.Range("J65").Value = 2990
$ws.Range("K65").Value = 96500
$ws.Range("L65").Value = 14950
$ws.Range("M65").Value = -93380
$ws.Range("N65").Value = -21190
$ws.Range("H116").Value = 4298.636
$ws.Range("I116").Value = 3022.5
$ws.Range("J116").Value = 4582.222
$ws.Range("K116").Value = 3022.5
$ws.Range("L116").Value = 4582.222
$ws.Range("M116").Value = 419.5
$ws.Range("N116").Value = -11466.222
$ws.Range("H132").Value = 4718.34
$ws.Range("I132").Value = 3252.55
$ws.Range("J132").Value = 9228.462
$ws.Range("K132").Value = 9757.650000000001
$ws.Range("L132").Value = 27685.386
$ws.Range("M132").Value = -7227.650000000001
$ws.Range("H135").Value = 424.21054
$ws.Range("I135").Value = 481.5625
$ws.Range("J135").Value = 118.333336
$ws.Range("K135").Value = 4334.0625
$ws.Range("L135").Value = 1065.000024
$ws.Range("M135").Value = -1799.0625
$ws.Range("N135").Value = -6135.000024
$ws.Range("H137").Value = 2562.9333
$ws.Range("I137").Value = 2445.8948
$ws.Range("J137").Value = 2765.0908
$ws.Range("K137").Value = 7337.6844
$ws.Range("L137").Value = 8295.2724
$ws.Range("M137").Value = -4787.6844
$ws.Range("H138").Value = 50002370
$ws.Range("I138").Value = 1665.7778
$ws.Range("J138").Value = 90912030
$ws.Range("K138").Value = 4997.3334
$ws.Range("L138").Value = 272736090
$ws.Range("M138").Value = 142.6665999999996
$ws.Range("N138").Value = -272746370
$ws.Range("H141").Value = 2006.0385
$ws.Range("I141").Value = 2095.0417
$ws.Range("J141").Value = 938
$ws.Range("K141").Value = 6285.125100000001
$ws.Range("L141").Value = 2814
$ws.Range("M141").Value = -1105.125100000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 110.25
$ws.Range("I4").Value = 84.59999999999999
$ws.Range("J4").Value = 153
$ws.Range("K4").Value = 84.59999999999999
$ws.Range("L4").Value = 153
$ws.Range("M4").Value = 31.40000000000001
$ws.Range("H32").Value = 3199.0508
$ws.Range("I32").Value = 2622.25
$ws.Range("J32").Value = 13966
$ws.Range("K32").Value = 2622.25
$ws.Range("L32").Value = 13966
$ws.Range("M32").Value = -2335.25
$ws.Range("H61").Value = 11372435
$ws.Range("I61").Value = 14711365
$ws.Range("J61").Value = 20069.9
$ws.Range("K61").Value = 14711365
$ws.Range("L61").Value = 20069.9
$ws.Range("M61").Value = -14711153
$ws.Range("N61").Value = -20493.9
$ws.Range("H74").Value = 4064.4285
$ws.Range("I74").Value = 3170.4583
$ws.Range("J74").Value = 6014.909
$ws.Range("K74").Value = 3170.4583
$ws.Range("L74").Value = 6014.909
$ws.Range("M74").Value = -2296.4583
$ws.Range("N74").Value = -7762.909
$ws.Range("H77").Value = 4064.4285
$ws.Range("I77").Value = 3170.4583
$ws.Range("J77").Value = 6014.909
$ws.Range("K77").Value = 15852.2915
$ws.Range("L77").Value = 30074.545
$ws.Range("M77").Value = -11484.2915
$ws.Range("N77").Value = -38810.545
$ws.Range("H110").Value = 5261.3213
$ws.Range("I110").Value = 4205.8096
$ws.Range("J110").Value = 8427.857
$ws.Range("K110").Value = 4205.8096
$ws.Range("L110").Value = 8427.857
$ws.Range("M110").Value = -2160.8096
$ws.Range("H136").Value = 11372435
$ws.Range("I136").Value = 14711365
$ws.Range("J136").Value = 20069.9
$ws.Range("K136").Value = 44134095
$ws.Range("L136").Value = 60209.7
$ws.Range("M136").Value = -44131545
$ws.Range("N136").Value = -65309.7

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 4779.0625
$ws.Range("I107").Value = 4506.364
$ws.Range("J107").Value = 5379
$ws.Range("K107").Value = 4506.364
$ws.Range("L107").Value = 5379
$ws.Range("M107").Value = -2586.364
$ws.Range("N107").Value = -9219
$ws.Range("H134").Value = 3277.3618
$ws.Range("I134").Value = 3424.5952
$ws.Range("J134").Value = 2040.6
$ws.Range("K134").Value = 10273.7856
$ws.Range("L134").Value = 6121.799999999999
$ws.Range("M134").Value = -7738.785600000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4054.9678
$ws.Range("I31").Value = 2069.6667
$ws.Range("J31").Value = 8224.1
$ws.Range("K31").Value = 2069.6667
$ws.Range("L31").Value = 8224.1
$ws.Range("M31").Value = -1774.6667
$ws.Range("N31").Value = -8814.1
$ws.Range("H34").Value = 4054.9678
$ws.Range("I34").Value = 2069.6667
$ws.Range("J34").Value = 8224.1
$ws.Range("K34").Value = 2069.6667
$ws.Range("L34").Value = 8224.1
$ws.Range("M34").Value = -1867.6667
$ws.Range("N34").Value = -8628.1
$ws.Range("H58").Value = 7358.7827
$ws.Range("I58").Value = 5941.091
$ws.Range("J58").Value = 8658.333000000001
$ws.Range("K58").Value = 5941.091
$ws.Range("L58").Value = 8658.333000000001
$ws.Range("M58").Value = -5738.091
$ws.Range("N58").Value = -9064.333000000001
$ws.Range("H80").Value = 56333.332
$ws.Range("I80").Value = 40000
$ws.Range("J80").Value = 64500
$ws.Range("K80").Value = 40000
$ws.Range("L80").Value = 64500
$ws.Range("M80").Value = -38877
$ws.Range("N80").Value = -66746
$ws.Range("H83").Value = 56333.332
$ws.Range("I83").Value = 40000
$ws.Range("J83").Value = 64500
$ws.Range("K83").Value = 120000
$ws.Range("L83").Value = 193500
$ws.Range("M83").Value = -114384
$ws.Range("N83").Value = -204732
$ws.Range("H122").Value = 2767.5186
$ws.Range("I122").Value = 2644.5
$ws.Range("J122").Value = 3013.5557
$ws.Range("K122").Value = 7933.5
$ws.Range("L122").Value = 9040.667099999999
$ws.Range("M122").Value = -5483.5
$ws.Range("N122").Value = -13940.6671
$ws.Range("H132").Value = 2728.12
$ws.Range("I132").Value = 2920.15
$ws.Range("J132").Value = 1960
$ws.Range("K132").Value = 8760.450000000001
$ws.Range("L132").Value = 5880
$ws.Range("M132").Value = -6230.450000000001
$ws.Range("H134").Value = 6126.2173
$ws.Range("I134").Value = 5375.625
$ws.Range("J134").Value = 7841.857
$ws.Range("K134").Value = 16126.875
$ws.Range("L134").Value = 23525.571
$ws.Range("M134").Value = -13591.875
$ws.Range("H136").Value = 7358.7827
$ws.Range("I136").Value = 5941.091
$ws.Range("J136").Value = 8658.333000000001
$ws.Range("K136").Value = 17823.273
$ws.Range("L136").Value = 25974.999
$ws.Range("M136").Value = -15273.273
$ws.Range("N136").Value = -31074.999
$ws.Range("H141").Value = 0
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").Value = $null

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H119").Value = 2499
$ws.Range("I119").Value = 2499
$ws.Range("J119").Value = 0
$ws.Range("K119").Value = 7497
$ws.Range("L119").Value = 0
$ws.Range("M119").Value = -2659
$ws.Range("H120").Value = 14151.214
$ws.Range("I120").Value = 3824.4
$ws.Range("J120").Value = 19888.334
$ws.Range("K120").Value = 11473.2
$ws.Range("L120").Value = 59665.00199999999
$ws.Range("M120").Value = -6635.200000000001
$ws.Range("N120").Value = -69341.00199999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H99").Value = 12523.667
$ws.Range("I99").Value = 8244.714
$ws.Range("J99").Value = 27500
$ws.Range("K99").Value = 8244.714
$ws.Range("L99").Value = 27500
$ws.Range("M99").Value = -5998.714
$ws.Range("N99").Value = -31992
$ws.Range("H102").Value = 3465
$ws.Range("I102").Value = 1670.8125
$ws.Range("J102").Value = 8249.5
$ws.Range("K102").Value = 1670.8125
$ws.Range("L102").Value = 8249.5
$ws.Range("M102").Value = -48.8125
$ws.Range("H113").Value = 288623.94
$ws.Range("I113").Value = 401989
$ws.Range("J113").Value = 5211.25
$ws.Range("K113").Value = 401989
$ws.Range("L113").Value = 5211.25
$ws.Range("M113").Value = -399819
$ws.Range("H132").Value = 1922.8148
$ws.Range("I132").Value = 1756.9166
$ws.Range("J132").Value = 3250
$ws.Range("K132").Value = 5270.7498
$ws.Range("L132").Value = 9750
$ws.Range("M132").Value = -2740.7498

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H96").Value = 43000
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 43000
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 43000
$ws.Range("N96").Value = -48492
$ws.Range("H132").Value = 11640.738
$ws.Range("I132").Value = 12568.581
$ws.Range("J132").Value = 9025.909
$ws.Range("K132").Value = 37705.743
$ws.Range("L132").Value = 27077.727
$ws.Range("M132").Value = -35175.743
$ws.Range("H136").Value = 4493.1763
$ws.Range("I136").Value = 4149.3438
$ws.Range("J136").Value = 9994.5
$ws.Range("K136").Value = 12448.0314
$ws.Range("L136").Value = 29983.5
$ws.Range("M136").Value = -9898.0314

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 1863292.1
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 1863292.1
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 1863292.1
$ws.Range("N5").Value = -1863516.1
$ws.Range("H41").Value = 17810.666
$ws.Range("I41").Value = 13000
$ws.Range("J41").Value = 18248
$ws.Range("K41").Value = 13000
$ws.Range("L41").Value = 18248
$ws.Range("M41").Value = -12610
$ws.Range("N41").Value = -19028
$ws.Range("H70").Value = 39554.43
$ws.Range("I70").Value = 33999
$ws.Range("J70").Value = 40480.332
$ws.Range("K70").Value = 33999
$ws.Range("L70").Value = 40480.332
$ws.Range("M70").Value = -33684
$ws.Range("N70").Value = -41110.332
$ws.Range("H73").Value = 39554.43
$ws.Range("I73").Value = 33999
$ws.Range("J73").Value = 40480.332
$ws.Range("K73").Value = 33999
$ws.Range("L73").Value = 40480.332
$ws.Range("M73").Value = -32907
$ws.Range("N73").Value = -42664.332
$ws.Range("H86").Value = 100000
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 100000
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 100000
$ws.Range("N86").Value = -102246
$ws.Range("H89").Value = 100000
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 100000
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 500000
$ws.Range("N89").Value = -511232
$ws.Range("H99").Value = 49999
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 49999
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 49999
$ws.Range("M99").Value = $null
$ws.Range("N99").Value = -55989
$ws.Range("H122").Value = 3423.524
$ws.Range("I122").Value = 1618.375
$ws.Range("J122").Value = 9200
$ws.Range("K122").Value = 4855.125
$ws.Range("L122").Value = 27600
$ws.Range("M122").Value = -2405.125
$ws.Range("H126").Value = 3082.0667
$ws.Range("I126").Value = 2636.5417
$ws.Range("J126").Value = 4864.1665
$ws.Range("K126").Value = 7909.625100000001
$ws.Range("L126").Value = 14592.4995
$ws.Range("M126").Value = -5439.625100000001
$ws.Range("N126").Value = -19532.4995
$ws.Range("H132").Value = 3243.5652
$ws.Range("I132").Value = 3243.5652
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 9730.695599999999
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -7200.695599999999
$ws.Range("H136").Value = 3758.3572
$ws.Range("I136").Value = 2868.3809
$ws.Range("J136").Value = 6428.2856
$ws.Range("K136").Value = 8605.1427
$ws.Range("L136").Value = 19284.8568
$ws.Range("M136").Value = -6055.1427
